# Daily attendance processing - swap the order of "Recorded By" entries
# so that the user email is listed before "System" wherever both
# appear together (column G: "System, dnasr281@gmail.com" ->
# "dnasr281@gmail.com, System"). Other column G values (e.g. "System"
# alone, or the email alone) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
